$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "Approved/Rejected" column (I) for rows 2-11 with "Approved",
# using the same text-formatted style already used for rows 12-29.
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Range("I$r")
    $cell.NumberFormat = "@"
    $cell.Value = "Approved"
}

# Match the author's saved view state: scrolled back to the top of the
# sheet (column I, row 1) with the newly-filled range selected.
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("I2:I11").Select()
